$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data block: replicate the formatting of the last existing
# data row (75) down across the ten new rows (76-85).
$ws.Range("A75:J75").Copy($ws.Range("A76:J85"))

# Row 81's "Tested (all)" cell was pasted with the plain/general number
# format (matching the rest of that row) instead of the thousands format
# used by the other rows in column B.
$ws.Range("C81").Copy()
$ws.Range("B81").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New daily COVID-19 figures for Slovenia, 2020-06-04 update.
$data = @(
  @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
  @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
  @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
  @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
  @(43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0),
  @(43981, 78793, 264, 1473, 0, 6, 2, 1, 108, 0),
  @(43982, 79039, 246, 1473, 0, 5, 1, 0, 109, 1),
  @(43983, 79698, 659, 1475, 2, 5, 1, 0, 109, 0),
  @(43984, 80505, 807, 1477, 2, 5, 0, 0, 109, 0),
  @(43985, 81333, 828, 1477, 0, 5, 0, 0, 109, 0)
)

$r = 76
foreach ($row in $data) {
  for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item($r, $c).Value2 = $row[$c - 1]
  }
  $r++
}

# Grow the worksheet table (ListObject) so the filter/banding covers the
# newly added rows too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J85"))

# Move the view down to the new bottom of the table, matching where a
# user would land after entering the last row of data.
$ws.Range("A85:J85").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1
